$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers (single dot)
# need to be forced to Text format first, otherwise Excel will
# auto-convert them to numeric values instead of keeping them as text.
$textCells = @("D5", "D6", "D11", "D14", "D15", "D16", "D19", "D20", "D25", "D27", "D39", "D42", "D45", "D47")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D5").Value = '210.25'
$ws.Range("D6").Value = '0.504'
$ws.Range("D11").Value = '0.0844'
$ws.Range("D14").Value = '4.07'
$ws.Range("D15").Value = '0.518'
$ws.Range("D16").Value = '64.37'
$ws.Range("D19").Value = '7.47'
$ws.Range("D20").Value = '210.77'
$ws.Range("D25").Value = '145.12'
$ws.Range("D27").Value = '7.05'
$ws.Range("D39").Value = '1.09'
$ws.Range("D42").Value = '5.62'
$ws.Range("D45").Value = '62.56'
$ws.Range("D47").Value = '87.93'

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}

# Remaining cell updates (already safe as text - contain letters, URLs,
# percent signs, or multiple dots so Excel keeps them as text)
$ws.Range("D2").Value = '26.359.59'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '1.590.04'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  -0.47%  '
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '1.816.76'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.612.93'
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '26.350.77'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("E19").Value = '  +5.18%  '
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("E23").Value = '  -3.76%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("E25").Value = '  +0.96%  '
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("E27").Value = '  -1.24%  '
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  -0.34%  '
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("D34").Value = '1.306.73'
$ws.Range("E34").Value = '  +2.11%  '
$ws.Range("E35").Value = '  +2.79%  '
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("E37").Value = '  -0.70%  '
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("E39").Value = '  -13.80%  '
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("E41").Value = '  -0.43%  '
$ws.Range("E42").Value = '  +3.78%  '
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").Value = '1.727.29'
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("E48").Value = '  -4.58%  '
$ws.Range("E49").Value = '  -4.46%  '
$ws.Range("E50").Value = '  -1.33%  '

